$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2271.8572
$ws.Range("I62").Value = 2271.8572
$ws.Range("K62").Value = 2271.8572
$ws.Range("M62").Value = -1647.8572
$ws.Range("H65").Value = 2271.8572
$ws.Range("I65").Value = 2271.8572
$ws.Range("K65").Value = 11359.286
$ws.Range("M65").Value = -8239.286
$ws.Range("H94").Value = 2309.25
$ws.Range("I94").Value = 2309.25
$ws.Range("K94").Value = 2309.25
$ws.Range("M94").Value = -1858.25
$ws.Range("H113").Value = 4169435.5
$ws.Range("I113").Value = 6669177
$ws.Range("K113").Value = 6669177
$ws.Range("M113").Value = -6665923
$ws.Range("H116").Value = 5915
$ws.Range("I116").Value = 6860.5557
$ws.Range("J116").Value = 3787.5
$ws.Range("K116").Value = 6860.5557
$ws.Range("L116").Value = 3787.5
$ws.Range("M116").Value = -3418.5557
$ws.Range("N116").Value = -10671.5
$ws.Range("H129").Value = 1067.65
$ws.Range("I129").Value = 761.4286
$ws.Range("J129").Value = 1108.0944
$ws.Range("K129").Value = 2284.2858
$ws.Range("L129").Value = 3324.2832
$ws.Range("M129").Value = 2715.7142
$ws.Range("N129").Value = -13324.2832
$ws.Range("H132").Value = 1400.3235
$ws.Range("I132").Value = 1158.862
$ws.Range("J132").Value = 2800.8
$ws.Range("K132").Value = 3476.586
$ws.Range("L132").Value = 8402.400000000001
$ws.Range("M132").Value = -946.5860000000002
$ws.Range("N132").Value = -13462.4
$ws.Range("H138").Value = 2032.8667
$ws.Range("I138").Value = 1202.909
$ws.Range("J138").Value = 2826.739
$ws.Range("K138").Value = 3608.727
$ws.Range("L138").Value = 8480.217000000001
$ws.Range("M138").Value = 1531.273
$ws.Range("N138").Value = -18760.217
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1214.4445
$ws.Range("I2").Value = 629.5
$ws.Range("J2").Value = 1682.4
$ws.Range("K2").Value = 629.5
$ws.Range("L2").Value = 1682.4
$ws.Range("M2").Value = -516.5
$ws.Range("N2").Value = -1908.4
$ws.Range("H32").Value = 8829.929
$ws.Range("I32").Value = 6943.25
$ws.Range("K32").Value = 6943.25
$ws.Range("M32").Value = -6656.25
$ws.Range("H61").Value = 6258.4
$ws.Range("I61").Value = 6678.1816
$ws.Range("J61").Value = 5104
$ws.Range("K61").Value = 6678.1816
$ws.Range("L61").Value = 5104
$ws.Range("M61").Value = -6466.1816
$ws.Range("N61").Value = -5528
$ws.Range("H74").Value = 1699.75
$ws.Range("I74").Value = 1516
$ws.Range("J74").Value = 2146
$ws.Range("K74").Value = 1516
$ws.Range("L74").Value = 2146
$ws.Range("M74").Value = -642
$ws.Range("N74").Value = -3894
$ws.Range("H77").Value = 1699.75
$ws.Range("I77").Value = 1516
$ws.Range("J77").Value = 2146
$ws.Range("K77").Value = 7580
$ws.Range("L77").Value = 10730
$ws.Range("M77").Value = -3212
$ws.Range("N77").Value = -19466
$ws.Range("H116").Value = 1214.4445
$ws.Range("I116").Value = 629.5
$ws.Range("J116").Value = 1682.4
$ws.Range("K116").Value = 629.5
$ws.Range("L116").Value = 1682.4
$ws.Range("M116").Value = 1664.5
$ws.Range("N116").Value = -6270.4
$ws.Range("H117").Value = 32300
$ws.Range("J117").Value = 32300
$ws.Range("L117").Value = 32300
$ws.Range("N117").Value = -41478
$ws.Range("H122").Value = 1510127.2
$ws.Range("I122").Value = 1833447.4
$ws.Range("K122").Value = 5500342.199999999
$ws.Range("M122").Value = -5497892.199999999
$ws.Range("H132").Value = 2043236.4
$ws.Range("I132").Value = 1528.8948
$ws.Range("J132").Value = 9096408
$ws.Range("K132").Value = 4586.6844
$ws.Range("L132").Value = 27289224
$ws.Range("M132").Value = -2056.6844
$ws.Range("N132").Value = -27294284
$ws.Range("H136").Value = 6258.4
$ws.Range("I136").Value = 6678.1816
$ws.Range("J136").Value = 5104
$ws.Range("K136").Value = 20034.5448
$ws.Range("L136").Value = 15312
$ws.Range("M136").Value = -17484.5448
$ws.Range("N136").Value = -20412
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1214.4445
$ws.Range("I3").Value = 629.5
$ws.Range("J3").Value = 1682.4
$ws.Range("K3").Value = 629.5
$ws.Range("L3").Value = 1682.4
$ws.Range("M3").Value = -515.5
$ws.Range("N3").Value = -1910.4
$ws.Range("H86").Value = 1774.6842
$ws.Range("I86").Value = 1735.9333
$ws.Range("J86").Value = 1920
$ws.Range("K86").Value = 1735.9333
$ws.Range("L86").Value = 1920
$ws.Range("M86").Value = -612.9332999999999
$ws.Range("N86").Value = -4166
$ws.Range("H89").Value = 1774.6842
$ws.Range("I89").Value = 1735.9333
$ws.Range("J89").Value = 1920
$ws.Range("K89").Value = 8679.666499999999
$ws.Range("L89").Value = 9600
$ws.Range("M89").Value = -3063.666499999999
$ws.Range("N89").Value = -20832
$ws.Range("H94").Value = 2674
$ws.Range("I94").Value = 2221.6
$ws.Range("J94").Value = 2997.1428
$ws.Range("K94").Value = 2221.6
$ws.Range("L94").Value = 2997.1428
$ws.Range("M94").Value = -1770.6
$ws.Range("N94").Value = -3899.1428
$ws.Range("H105").Value = 2908.2632
$ws.Range("I105").Value = 2149.8333
$ws.Range("J105").Value = 3258.3076
$ws.Range("K105").Value = 2149.8333
$ws.Range("L105").Value = 3258.3076
$ws.Range("M105").Value = -402.8332999999998
$ws.Range("N105").Value = -6752.3076
$ws.Range("H134").Value = 5044.5586
$ws.Range("I134").Value = 5689.0386
$ws.Range("J134").Value = 2950
$ws.Range("K134").Value = 17067.1158
$ws.Range("L134").Value = 8850
$ws.Range("M134").Value = -14532.1158
$ws.Range("N134").Value = -13920
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3374.9583
$ws.Range("I99").Value = 3642.8096
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 3642.8096
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -2144.8096
$ws.Range("N99").Value = -4496
$ws.Range("H122").Value = 2528801.2
$ws.Range("I122").Value = 5556151.5
$ws.Range("K122").Value = 16668454.5
$ws.Range("M122").Value = -16666004.5
$ws.Range("H126").Value = 3374.9583
$ws.Range("I126").Value = 3642.8096
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 10928.4288
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -8458.4288
$ws.Range("N126").Value = -9440
$ws.Range("H132").Value = 2137.5
$ws.Range("I132").Value = 1733.45
$ws.Range("K132").Value = 5200.35
$ws.Range("M132").Value = -2670.35
$ws.Range("H134").Value = 2401.913
$ws.Range("I134").Value = 2461.6052
$ws.Range("J134").Value = 2118.375
$ws.Range("K134").Value = 7384.8156
$ws.Range("L134").Value = 6355.125
$ws.Range("M134").Value = -4849.8156
$ws.Range("N134").Value = -11425.125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H122").Value = 2118.951
$ws.Range("I122").Value = 266.0909
$ws.Range("J122").Value = 2526.58
$ws.Range("K122").Value = 2394.8181
$ws.Range("L122").Value = 22739.22
$ws.Range("M122").Value = 55.18190000000004
$ws.Range("N122").Value = -27639.22
$ws.Range("H129").Value = 716.125
$ws.Range("I129").Value = 675.5714
$ws.Range("K129").Value = 2026.7142
$ws.Range("M129").Value = 2973.2858
$ws.Range("H136").Value = 4441.2
$ws.Range("I136").Value = 11895.556
$ws.Range("J136").Value = 2804.878
$ws.Range("K136").Value = 35686.66800000001
$ws.Range("L136").Value = 8414.634
$ws.Range("M136").Value = -30586.66800000001
$ws.Range("N136").Value = -18614.634
$ws.Range("H137").Value = 6797.9165
$ws.Range("I137").Value = 5410.8
$ws.Range("J137").Value = 8305.652
$ws.Range("K137").Value = 16232.4
$ws.Range("L137").Value = 24916.956
$ws.Range("M137").Value = -11132.4
$ws.Range("N137").Value = -35116.956
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1413987.4
$ws.Range("I102").Value = 2422307
$ws.Range("J102").Value = 2340
$ws.Range("K102").Value = 2422307
$ws.Range("L102").Value = 2340
$ws.Range("M102").Value = -2420685
$ws.Range("N102").Value = -5584
$ws.Range("H132").Value = 2279.1538
$ws.Range("I132").Value = 1739.5385
$ws.Range("J132").Value = 2818.7693
$ws.Range("K132").Value = 5218.6155
$ws.Range("L132").Value = 8456.3079
$ws.Range("M132").Value = -2688.6155
$ws.Range("N132").Value = -13516.3079
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12351400
$ws.Range("I132").Value = 18525306
$ws.Range("J132").Value = 3586.4443
$ws.Range("K132").Value = 55575918
$ws.Range("L132").Value = 10759.3329
$ws.Range("M132").Value = -55573388
$ws.Range("N132").Value = -15819.3329
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 968.2
$ws.Range("I113").Value = 861.8570999999999
$ws.Range("J113").Value = 1216.3334
$ws.Range("K113").Value = 2585.5713
$ws.Range("L113").Value = 3649.0002
$ws.Range("M113").Value = -415.5712999999996
$ws.Range("N113").Value = -7989.0002
$ws.Range("H132").Value = 1676.5128
$ws.Range("I132").Value = 1359.45
$ws.Range("J132").Value = 2010.2632
$ws.Range("K132").Value = 4078.35
$ws.Range("L132").Value = 6030.7896
$ws.Range("M132").Value = -1548.35
$ws.Range("N132").Value = -11090.7896
$ws.Range("H136").Value = 1782.238
$ws.Range("I136").Value = 1770.4324
$ws.Range("J136").Value = 1799.0385
$ws.Range("K136").Value = 5311.2972
$ws.Range("L136").Value = 5397.1155
$ws.Range("M136").Value = -2761.2972
$ws.Range("N136").Value = -10497.1155

Write-Host "Applied 247 cell changes"